# "datatill 18 May 12PM"
# Correct the 2021-05-08 Ordered Amount on the Wallet ledger, append the
# new ledger rows through 2021-05-14, and record the new
# Device-Connectivity checkpoint row for 2021-05-18.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Wallet": fix row 206 and append the new ledger rows 207-215
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Wallet")

# Row 206 (2021-05-08) Ordered Amount was mis-entered; fix it. The
# Wallet Balance formula in E206 recalculates automatically.
$ws.Range("B206").Value = 33264

# New transaction rows, in order. Debit -> "Ordered Amount" column B;
# Credit -> "Manual Added" column C (mirrors the existing rows above).
$newRows = @(
    @{ Row = 207; Date = 44325; Debit = 34304;  Credit = $null   },
    @{ Row = 208; Date = 44326; Debit = 42620;  Credit = $null   },
    @{ Row = 209; Date = 44327; Debit = 31185;  Credit = $null   },
    @{ Row = 210; Date = 44328; Debit = 38462;  Credit = $null   },
    @{ Row = 211; Date = 44329; Debit = 33264;  Credit = $null   },
    @{ Row = 212; Date = 44329; Debit = $null;  Credit = 16745   },
    @{ Row = 213; Date = 44329; Debit = $null;  Credit = 100237  },
    @{ Row = 214; Date = 44330; Debit = 36383;  Credit = $null   }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $prev = $r - 1

    # Clone formatting from the matching template row above (date style
    # for column A, "debit"/"credit" mode style for column D) before
    # writing the new values.
    $ws.Range("A206").Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)

    if ($item.Debit -ne $null) {
        $ws.Range("D206").Copy()
    } else {
        $ws.Range("D205").Copy()
    }
    $ws.Cells.Item($r, 4).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $item.Date

    if ($item.Debit -ne $null) {
        $ws.Cells.Item($r, 2).Value = $item.Debit
        $ws.Cells.Item($r, 4).Value = "Ordered Amount"
    } else {
        $ws.Cells.Item($r, 3).Value = $item.Credit
        $ws.Cells.Item($r, 4).Value = "Manual Added"
    }

    $ws.Cells.Item($r, 5).Formula = '=IF(A' + $r + '="","",SUM(E' + $prev + '-B' + $r + '+C' + $r + '))'
}

# Trailing blank row, keeping the date-column styling used throughout.
$ws.Range("A206").Copy()
$ws.Range("A215").PasteSpecial(-4122)

[void]$ws.Range("B214").Select()

# ---------------------------------------------------------------------
# "Device-Connectivity": append the feed row for the new "data till"
# checkpoint date (2021-05-18).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Device-Connectivity")
$ws2.Range("A7:B7").Copy()
$ws2.Range("A8:B8").PasteSpecial(-4122)
$ws2.Range("A8").Value = 25451
$ws2.Range("B8").Value = 44334

[void]$ws2.Range("G11").Select()
